# Update "想去人数" (wish-to-go count) figures in sheet "展览" (sheet1)
# and the mirrored rows in sheet "全部类型" (sheet4) to match the
# newly generated output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws1.Range("F2").Value  = 14927
$ws1.Range("F3").Value  = 18696
$ws1.Range("F16").Value = 58
$ws1.Range("F17").Value = 1432
$ws1.Range("F22").Value = 7775
$ws1.Range("F24").Value = 25
$ws1.Range("F28").Value = 5980
$ws1.Range("F29").Value = 109
$ws1.Range("F31").Value = 163
$ws1.Range("F33").Value = 265
$ws1.Range("F34").Value = 5351

$ws4 = $wb.Worksheets.Item(4)   # 全部类型
$ws4.Range("F2").Value  = 14927
$ws4.Range("F3").Value  = 18696
$ws4.Range("F16").Value = 58
$ws4.Range("F17").Value = 1432
$ws4.Range("F23").Value = 7775
$ws4.Range("F25").Value = 25
$ws4.Range("F31").Value = 5980
$ws4.Range("F32").Value = 109
$ws4.Range("F34").Value = 163
$ws4.Range("F36").Value = 265
$ws4.Range("F37").Value = 5351
